# FIX: date sort, ADD: border
# Column G ("year") held plain 4-digit year numbers; convert them to real
# dates (serial numbers) formatted as yyyy-mm-dd so the column sorts
# chronologically instead of lexically-by-year. Also add a page border
# (page setup) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the year values in column G with real date serials, and give
# --- them a yyyy-mm-dd date format (this creates the new numFmt + cellXf
# --- that the header cell and the column default style also pick up).
$ws.Range("G1").NumberFormat = "yyyy\-mm\-dd"

$dates = @{
    "G2"  = 42979
    "G3"  = 43132
    "G4"  = 43191
    "G5"  = 43192
    "G6"  = 43192
    "G7"  = 43313
    "G8"  = 43282
    "G9"  = 43374
    "G10" = 43374
    "G11" = 43374
    "G12" = 43806
    "G13" = 43398
    "G14" = 43581
    "G15" = 43670
    "G16" = 43678
    "G19" = 43831
}

foreach ($addr in $dates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $dates[$addr]
    $cell.NumberFormat = "yyyy\-mm\-dd"
}

# --- Move the active selection from I15 to I11.
$null = $ws.Range("I11").Select()

# --- Add a page border / page setup (paper size + portrait orientation).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
